# LogicComponent.pptx update (Issue 724 - update DevMan diagrams to match
# design changes).
#
# This slide gains a new "GateKeeper" box feeding into the existing
# AccountsLogic / CoursesLogic / EvaluationsLogic stack, so that stack of
# rectangles and its connectors slide down / shrink to make room.

function EmuToPt($targetEmu) {
    # Shape.Left/Top/Width/Height round-trip through a Single (32-bit
    # float) internally, same as real PowerPoint, so a naive EMU/12700
    # division can land 1 EMU off after it is cast back down to Single.
    # Nudge the point value in tiny steps until it lands back on the
    # exact requested EMU value once re-quantized to Single precision.
    $pts = $targetEmu / 12700.0
    for ($i = 0; $i -lt 4000; $i++) {
        $f32 = [double]([float]$pts)
        $emu = [math]::Floor($f32 * 12700 + 0.0000001)
        if ($emu -eq $targetEmu) {
            return $pts
        }
        if ($emu -lt $targetEmu) {
            $pts += 0.0000005
        } else {
            $pts -= 0.0000005
        }
    }
    return $targetEmu / 12700.0
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ------------------------------------------------------------------
# 1) Add the new "GateKeeper" rectangle.
#    Clone "Rectangle 39" (AccountsLogic) *before* it is repositioned
#    below, so the new shape inherits the exact same style/formatting
#    (accent4 fill/line/font refs, 1600/bold run formatting, etc.)
# ------------------------------------------------------------------
$accountsLogic = $s.Shapes.Item("Rectangle 39")
$gateKeeper = $accountsLogic.Duplicate()
$gateKeeper.Name = "Rectangle 54"
$gateKeeper.Left = EmuToPt 5486400
$gateKeeper.Top = EmuToPt 1295400
$gateKeeper.Width = EmuToPt 1676400
$gateKeeper.Height = EmuToPt 457200
$gateKeeper.TextFrame.TextRange.Text = "GateKeeper"

# ------------------------------------------------------------------
# 2) Add the new connector running out of the GateKeeper box.
#    Clone "Straight Arrow Connector 81" (which already carries the
#    right-middle connection-site + dashed/arrow line formatting)
#    before it is repositioned below.
# ------------------------------------------------------------------
$templateConnector = $s.Shapes.Item("Straight Arrow Connector 81")
$gateKeeperConnector = $templateConnector.Duplicate()
$gateKeeperConnector.Name = "Straight Arrow Connector 58"
$gateKeeperConnector.Left = EmuToPt 7162800
$gateKeeperConnector.Top = EmuToPt 1524000
$gateKeeperConnector.Width = EmuToPt 762000
$gateKeeperConnector.Height = 0

# ------------------------------------------------------------------
# 3) Reposition/resize the pre-existing shapes to make room.
# ------------------------------------------------------------------

# TeamEvalResult box shifts down.
$teamEvalResult = $s.Shapes.Item("Rectangle 85")
$teamEvalResult.Top = EmuToPt 3581400

# Emails box shifts down.
$emails = $s.Shapes.Item("Rectangle 115")
$emails.Top = EmuToPt 4040088

# Connector feeding into the Emails box shifts down with it.
$emailsConnector = $s.Shapes.Item("Straight Arrow Connector 47")
$emailsConnector.Top = EmuToPt 4191744

# AccountsLogic moves down and shrinks to make room for GateKeeper above it.
$accountsLogic.Top = EmuToPt 1828800
$accountsLogic.Height = EmuToPt 457200

# CoursesLogic moves down and shrinks.
$coursesLogic = $s.Shapes.Item("Rectangle 55")
$coursesLogic.Top = EmuToPt 2362200
$coursesLogic.Height = EmuToPt 381000

# EvaluationsLogic moves down and shrinks.
$evaluationsLogic = $s.Shapes.Item("Rectangle 57")
$evaluationsLogic.Top = EmuToPt 2819400
$evaluationsLogic.Height = EmuToPt 381000

# The three dashed connectors on the right edge move down to track
# CoursesLogic / EvaluationsLogic / AccountsLogic respectively.
$connector67 = $s.Shapes.Item("Straight Arrow Connector 67")
$connector67.Top = EmuToPt 2514600

$connector79 = $s.Shapes.Item("Straight Arrow Connector 79")
$connector79.Top = EmuToPt 3048000

$connector81 = $s.Shapes.Item("Straight Arrow Connector 81")
$connector81.Top = EmuToPt 2057400
